$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 14568844
$ws.Range("I2").Value = 11904831
$ws.Range("J2").Value = 19230868
$ws.Range("K2").Value = 11904831
$ws.Range("L2").Value = 19230868
$ws.Range("M2").Value = -11904718
$ws.Range("N2").Value = -19231094
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H125").Value = 1161.0834
$ws.Range("I125").Value = 900
$ws.Range("J125").Value = 1213.3
$ws.Range("K125").Value = 8100
$ws.Range("L125").Value = 10919.7
$ws.Range("M125").Value = -5640
$ws.Range("N125").Value = -15839.7
$ws.Range("H127").Value = 1164.0769
$ws.Range("I127").Value = 654.1111
$ws.Range("J127").Value = 1434.0588
$ws.Range("K127").Value = 1962.3333
$ws.Range("L127").Value = 4302.1764
$ws.Range("M127").Value = 2997.6667
$ws.Range("N127").Value = -14222.1764
$ws.Range("H129").Value = 996.4286
$ws.Range("I129").Value = 1411.8334
$ws.Range("J129").Value = 952.7018
$ws.Range("K129").Value = 4235.5002
$ws.Range("L129").Value = 2858.1054
$ws.Range("M129").Value = 764.4997999999996
$ws.Range("N129").Value = -12858.1054
$ws.Range("H132").Value = 11635007
$ws.Range("I132").Value = 13519917
$ws.Range("J132").Value = 11390.833
$ws.Range("K132").Value = 40559751
$ws.Range("L132").Value = 34172.499
$ws.Range("M132").Value = -40557221
$ws.Range("N132").Value = -39232.499

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10417524
$ws.Range("I74").Value = 11111979
$ws.Range("J74").Value = 693
$ws.Range("K74").Value = 11111979
$ws.Range("L74").Value = 693
$ws.Range("M74").Value = -11111105
$ws.Range("N74").Value = -2441
$ws.Range("H77").Value = 10417524
$ws.Range("I77").Value = 11111979
$ws.Range("J77").Value = 693
$ws.Range("K77").Value = 55559895
$ws.Range("L77").Value = 3465
$ws.Range("M77").Value = -55555527
$ws.Range("N77").Value = -12201
$ws.Range("H80").Value = 27110
$ws.Range("J80").Value = 27110
$ws.Range("L80").Value = 27110
$ws.Range("N80").Value = -29106
$ws.Range("H83").Value = 27110
$ws.Range("J83").Value = 27110
$ws.Range("L83").Value = 81330
$ws.Range("N83").Value = -91314
$ws.Range("H102").Value = 1548.5714
$ws.Range("I102").Value = 1368
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1368
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 254
$ws.Range("N102").Value = -5244

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 535.4545000000001
$ws.Range("I94").Value = 537.4
$ws.Range("J94").Value = 531.2857
$ws.Range("K94").Value = 537.4
$ws.Range("L94").Value = 531.2857
$ws.Range("M94").Value = -86.39999999999998
$ws.Range("N94").Value = -1433.2857
$ws.Range("H99").Value = 1681.8518
$ws.Range("I99").Value = 1016.8461
$ws.Range("K99").Value = 1016.8461
$ws.Range("M99").Value = 481.1539

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 466.33334
$ws.Range("I19").Value = 466.33334
$ws.Range("K19").Value = 466.33334
$ws.Range("M19").Value = -296.33334
$ws.Range("H24").Value = 466.33334
$ws.Range("I24").Value = 466.33334
$ws.Range("K24").Value = 466.33334
$ws.Range("M24").Value = -296.33334
$ws.Range("H31").Value = 3751.518
$ws.Range("I31").Value = 2889.566
$ws.Range("J31").Value = 5274.3
$ws.Range("K31").Value = 2889.566
$ws.Range("L31").Value = 5274.3
$ws.Range("M31").Value = -2594.566
$ws.Range("N31").Value = -5864.3
$ws.Range("H34").Value = 3751.518
$ws.Range("I34").Value = 2889.566
$ws.Range("J34").Value = 5274.3
$ws.Range("K34").Value = 2889.566
$ws.Range("L34").Value = 5274.3
$ws.Range("M34").Value = -2687.566
$ws.Range("N34").Value = -5678.3
$ws.Range("H122").Value = 16668233
$ws.Range("I122").Value = 22727954
$ws.Range("K122").Value = 68183862
$ws.Range("M122").Value = -68181412

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5884459.5
$ws.Range("I80").Value = 2253.6
$ws.Range("J80").Value = 50001004
$ws.Range("K80").Value = 2253.6
$ws.Range("L80").Value = 50001004
$ws.Range("M80").Value = -1255.6
$ws.Range("N80").Value = -50003000
$ws.Range("H83").Value = 5884459.5
$ws.Range("I83").Value = 2253.6
$ws.Range("J83").Value = 50001004
$ws.Range("K83").Value = 11268
$ws.Range("L83").Value = 250005020
$ws.Range("M83").Value = -6276
$ws.Range("N83").Value = -250015004
$ws.Range("H102").Value = 1048.6471
$ws.Range("J102").Value = 1037.6666
$ws.Range("L102").Value = 1037.6666
$ws.Range("N102").Value = -4281.6666
$ws.Range("H122").Value = 22734738
$ws.Range("I122").Value = 27785818
$ws.Range("J122").Value = 4875.25
$ws.Range("K122").Value = 83357454
$ws.Range("L122").Value = 14625.75
$ws.Range("M122").Value = -83355004
$ws.Range("N122").Value = -19525.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 733
$ws.Range("I30").Value = 733
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 733
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -625
$ws.Range("N30").ClearContents()
$ws.Range("H100").Value = 3126.5454
$ws.Range("I100").Value = 2389.5454
$ws.Range("J100").Value = 3863.5454
$ws.Range("K100").Value = 2389.5454
$ws.Range("L100").Value = 3863.5454
$ws.Range("M100").Value = -1848.5454
$ws.Range("N100").Value = -4945.5454
$ws.Range("H122").Value = 6181.615
$ws.Range("I122").Value = 6827.5713
$ws.Range("J122").Value = 3468.6
$ws.Range("K122").Value = 20482.7139
$ws.Range("L122").Value = 10405.8
$ws.Range("M122").Value = -18032.7139
$ws.Range("N122").Value = -15305.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H41").Value = 5938.5
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 5938.5
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 5938.5
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -6718.5
$ws.Range("H81").Value = 35714980
$ws.Range("J81").Value = 913.5
$ws.Range("L81").Value = 1827
$ws.Range("N81").Value = -3949
$ws.Range("H84").Value = 35714980
$ws.Range("J84").Value = 913.5
$ws.Range("L84").Value = 9135
$ws.Range("N84").Value = -19743
$ws.Range("H122").Value = 2021.3143
$ws.Range("I122").Value = 1545.3334
$ws.Range("J122").Value = 3059.818
$ws.Range("K122").Value = 4636.0002
$ws.Range("L122").Value = 9179.454000000002
$ws.Range("M122").Value = -2186.0002
$ws.Range("N122").Value = -14079.454
$ws.Range("H126").Value = 1765.4546
$ws.Range("I126").Value = 802.8570999999999
$ws.Range("J126").Value = 3450
$ws.Range("K126").Value = 2408.5713
$ws.Range("L126").Value = 10350
$ws.Range("M126").Value = 61.42870000000039
$ws.Range("N126").Value = -15290
$ws.Range("H132").Value = 25772.13
$ws.Range("I132").Value = 42145.4
$ws.Range("J132").Value = 6280.143
$ws.Range("K132").Value = 126436.2
$ws.Range("L132").Value = 18840.429
$ws.Range("M132").Value = -123906.2
$ws.Range("N132").Value = -23900.429
$ws.Range("H136").Value = 3216.8867
$ws.Range("I136").Value = 4593.148
$ws.Range("J136").Value = 1787.6923
$ws.Range("K136").Value = 13779.444
$ws.Range("L136").Value = 5363.0769
$ws.Range("M136").Value = -11229.444
$ws.Range("N136").Value = -10463.0769
